$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-22 Wednesday" "2025-01-23 Thursday"

Replace-Text "961×2=" "153×7="
Replace-Text "744×9=" "506×9="
Replace-Text "785×7=" "993×7="
Replace-Text "891×8=" "124×6="
Replace-Text "313×7=" "220×4="
Replace-Text "658×4=" "517×8="
Replace-Text "632×4=" "235×4="
Replace-Text "881×5=" "239×9="
Replace-Text "764×8=" "668×4="
Replace-Text "756×4=" "758×7="
Replace-Text "665×3=" "800×4="
Replace-Text "481×6=" "141×7="
Replace-Text "607×7=" "971×2="
Replace-Text "144×7=" "763×9="
Replace-Text "552×5=" "741×5="
Replace-Text "751×3=" "364×8="
Replace-Text "918×6=" "852×6="
Replace-Text "181×7=" "253×7="
Replace-Text "680×9=" "991×5="
Replace-Text "421×7=" "532×9="
Replace-Text "720×2=" "193×3="
Replace-Text "431×6=" "493×7="
Replace-Text "422×8=" "808×5="
Replace-Text "684×9=" "610×8="
Replace-Text "387×9=" "640×3="
